$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 125
$ws.Range("B1").Value = 200
$ws.Range("C1").Value = 125
$ws.Range("A2").Value = 125
$ws.Range("B2").Value = 179
$ws.Range("C2").Value = 125
